# auto increment public id and fix unit tests
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Parent Public ID column (A): rows 2-4 belong to "statement-01",
# rows 5-6 belong to "statement-02" (renamed from "Statement 1"/"Statement 6"
# and the group boundary moved down one row).
$ws.Range("A2").Value = "statement-01"
$ws.Range("A3").Value = "statement-01"
$ws.Range("A4").Value = "statement-01"
$ws.Range("A5").Value = "statement-02"
$ws.Range("A6").Value = "statement-02"

# Public ID column (B) is now auto incremented/generated, so the stored
# sample values are cleared out (style is preserved).
$ws.Range("B2").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()
